# Applies the "add logs for results examples" commit:
#  - adds a new "Extend:" shared-string label used at C8 and I37
#  - appends a 5th QRE data point (row 8) to the first results table
#  - appends a 4th qlog data point (row 37) to the second results table
#  - updates both scatter charts so the newly extended ranges are
#    plotted with circular markers, axis titles, and a bottom legend
#  - moves the sheet view to show the newly added rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Worksheet data — extend the two results tables with one more row
# ---------------------------------------------------------------------

# First table (rows 3-29): add row 8, mirroring the "Extend:" label used
# later in the second table at I37.
$ws.Range("C8").Value = "Extend:"
$ws.Range("E8").Value = 142.53700000000001
$ws.Range("G8").Value = 818.59

# Second table (rows 33-41): add a 4th data row under row 37.
$ws.Range("I37").Value = "Extend:"
$ws.Range("J37").Value = 125
$ws.Range("K37").Value = 764.37599999999998

# ---------------------------------------------------------------------
# 2. Chart 1 (first results table chart)
# ---------------------------------------------------------------------
$chart1 = $ws.ChartObjects().Item(1).Chart

$s1 = $chart1.SeriesCollection().Item(1)
$s1.XValues = $ws.Range("G4:G8")
$s1.Values = $ws.Range("E4:E8")
$s1.MarkerStyle = 8
$s1.MarkerSize = 5

$s2 = $chart1.SeriesCollection().Item(2)
$s2.MarkerStyle = 8
$s2.MarkerSize = 5

$chart1.Axes(1).HasTitle = $true
$chart1.Axes(1).AxisTitle.Text = "Time since connection started (ms)"
$chart1.Axes(2).HasTitle = $true
$chart1.Axes(2).AxisTitle.Text = "RTT (ms)"

$chart1.Legend.Position = -4107

# ---------------------------------------------------------------------
# 3. Chart 2 (second results table chart)
# ---------------------------------------------------------------------
$chart2 = $ws.ChartObjects().Item(2).Chart

$s3 = $chart2.SeriesCollection().Item(2)
$s3.XValues = $ws.Range("K34:K37")
$s3.Values = $ws.Range("J34:J37")

$chart2.Axes(1).HasTitle = $true
$chart2.Axes(1).AxisTitle.Text = "Time since connection started (ms)"
$chart2.Axes(2).HasTitle = $true
$chart2.Axes(2).AxisTitle.Text = "RTT (ms)"

$chart2.Legend.Position = -4107

# ---------------------------------------------------------------------
# 4. Sheet view — scroll down to show the newly extended rows
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("S21").Select()
